# Weekly fruit/vegetable price update: a new "Ají" (Inferno, Primera) record
# for 2022-09-23 is inserted as a new data row right before the existing
# row 64, pushing all subsequent rows (old 64..102) down by one (new
# 65..103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64 (shifts rows 64-102 down to 65-103).
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 44827
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112021
$ws.Cells.Item(64, 7).Value = "Ají"
$ws.Cells.Item(64, 8).Value = "Inferno"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 130
$ws.Cells.Item(64, 11).Value = 15000
$ws.Cells.Item(64, 12).Value = 16000
$ws.Cells.Item(64, 13).Value = 15500
$ws.Cells.Item(64, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 1033
$ws.Cells.Item(64, 17).Value = 15
$ws.Cells.Item(64, 18).Value = "Hortaliza"
